$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 4 with the latest test-mail entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A4").Value = "Kun jij dit even regelen?"
$ws.Range("B4").Value = "mailmind.test@zohomail.eu"
$ws.Range("C4").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("D4").Value = "Overig"
$ws.Range("E4").Value = "Beste klant,`nBedankt voor uw e-mail. Om u beter van dienst te kunnen zijn, zou u wat meer specifieke informatie kunnen verstrekken over wat u precies geregeld wilt hebben?`nMet vriendelijke groet,`n[Jouw naam]`nE-mailassistent"
$ws.Range("F4").Value = "2025-08-01 23:01:32"
$ws.Range("G4").Value = "Ja"
$ws.Range("H4").Value = "Nee"
$ws.Range("I4").Value = "Ja"
$ws.Range("J4").Value = "Nee"

# Extend the conditional-formatting ranges from rows 2:3 to 2:4
$ws.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D4"))
$ws.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G4"))
$ws.Range("H2:H3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H4"))
$ws.Range("I2:I3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I4"))
$ws.Range("J2:J3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J4"))

# --- Dashboard sheet: bump the "Overig" tally to reflect the new row ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 3
